$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New D/E values (dic_nbre_clients_poisson_2_keys / dic_nbre_clients_prob_poisson_2_values)
# per former row index (A column value 0..50), rows 2..52 after the edit.
$data = @(
    @(0, 0.153),
    @(2, 0.001),
    @(3, 0.007),
    @(4, 0.011),
    @(5, 0.02),
    @(6, 0.025),
    @(7, 0.039),
    @(8, 0.055),
    @(9, 0.048),
    @(10, 0.025),
    @(11, 0.03),
    @(12, 0.028),
    @(13, 0.026),
    @(14, 0.028),
    @(15, 0.048),
    @(16, 0.029),
    @(17, 0.036),
    @(18, 0.033),
    @(19, 0.026),
    @(20, 0.027),
    @(21, 0.017),
    @(22, 0.026),
    @(23, 0.031),
    @(24, 0.022),
    @(25, 0.018),
    @(26, 0.026),
    @(27, 0.012),
    @(28, 0.02),
    @(29, 0.01),
    @(30, 0.016),
    @(31, 0.011),
    @(32, 0.011),
    @(33, 0.015),
    @(34, 0.007),
    @(35, 0.011),
    @(36, 0.007),
    @(37, 0.009000000000000001),
    @(38, 0.009000000000000001),
    @(39, 0.005),
    @(40, 0.001),
    @(41, 0.003),
    @(43, 0.005),
    @(45, 0.001),
    @(46, 0.002),
    @(47, 0.003),
    @(48, 0.001),
    @(49, 0.001),
    @(50, 0.001),
    @(53, 0.001),
    @(54, 0.001),
    @(57, 0.001)
)

# First, delete the two trailing rows (53 and 54) so the used range shrinks
# from A1:E54 to A1:E52.
$ws.Rows.Item(54).Delete()
$ws.Rows.Item(53).Delete()

# Update B and C columns (lamda_1 / lamda_2) for every data row, 2..52.
$ws.Range("B2:B52").Value = 33.94444444444444
$ws.Range("C2:C52").Value = 1.95

# Update D and E columns (keys / probabilities) row by row.
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 4).Value = $data[$i][0]
    $ws.Cells.Item($r, 5).Value = $data[$i][1]
}
